$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, pushing existing rows 8-20 down to 9-21.
$ws.Rows(8).Insert()

# Populate the new row 8 with the latest weekly price entry.
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = 44935
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = 100112006
$ws.Cells.Item(8, 7).Value = "Repollo"
$ws.Cells.Item(8, 8).Value = "Copenhague"
$ws.Cells.Item(8, 9).Value = "Segunda"
$ws.Cells.Item(8, 10).Value = 1000
$ws.Cells.Item(8, 11).Value = 400
$ws.Cells.Item(8, 12).Value = 500
$ws.Cells.Item(8, 13).Value = 460
$ws.Cells.Item(8, 14).Value = "$/unidad"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 460
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = "Hortaliza"
